# Add a new row of HICP data for 11/2022 (row 13), extend the cumulative
# formulas down into it, restyle the previously-blank label cells below it
# (A14:A29) to match the "text" style already used by A12/A13, and move the
# active selection/cell shown when the sheet is opened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data (row 13: "11/2022") -----------------------------
$ws.Range("A13").Value = "11/2022"
$ws.Range("B13").Value = 1.3
$ws.Range("C13").Formula = "=B13/100"
$ws.Range("D13").Formula = "=(C13+1)*D12"

# Match the text-number-format style already used for the month labels
# (A12 uses it; it's the style that was introduced for the last label row).
$ws.Range("A12").Copy()
$ws.Range("A13:A29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Move the active selection shown on open --------------------------
$ws.Range("H19").Select()
